$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 23:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 815491
$ws.Range("C4").Value = 22732
$ws.Range("E4").Value = 687774
$ws.Range("G4").Value = 2583
$ws.Range("H4").Value = 45097

# Row 8 - Alemania
$ws.Range("B8").Value = 148291
$ws.Range("C8").Value = 1226
$ws.Range("E8").Value = 48058
$ws.Range("G8").Value = 171
$ws.Range("H8").Value = 5033

# Row 14 - Brasil
$ws.Range("B14").Value = 43079
$ws.Range("C14").Value = 2336
$ws.Range("E14").Value = 17347
$ws.Range("G14").Value = 154
$ws.Range("H14").Value = 2741

# Row 16 - Canada
$ws.Range("B16").Value = 38210
$ws.Range("C16").Value = 1381
$ws.Range("D16").Value = 13143
$ws.Range("E16").Value = 23236

# Row 151 - Barbados
$ws.Range("D151").Value = 25
$ws.Range("E151").Value = 45

# Row 202 - Nicaragua
$ws.Range("D202").Value = 7
$ws.Range("E202").Value = 1
